# Summer 2024 Working Hours.xlsx - add "Day of week" column + new day's entry
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before the old "Hours" column (B) to hold the day of
# week abbreviation. This shifts Hours/spacer/Notes/Links/TODO one column to
# the right (B->C, C->D, D->E, E->F, F->G) and carries over the date-style
# formatting already applied to column A for every existing row.
$ws.Columns("B").Insert()

# Approximate the narrow width used for the new day-of-week column.
$ws.Columns("B").ColumnWidth = 1.86

# Header
$ws.Range("B1").Value = "Day of week"

# Day-of-week abbreviations for each existing dated row (Monday..Sunday ->
# M / T / W / T / F / S / S)
$ws.Range("B12").Value = "M"
$ws.Range("B13").Value = "T"
$ws.Range("B14").Value = "W"
$ws.Range("B15").Value = "T"
$ws.Range("B16").Value = "F"
$ws.Range("B17").Value = "S"
$ws.Range("B18").Value = "S"
$ws.Range("B19").Value = "T"
$ws.Range("B20").Value = "W"
$ws.Range("B21").Value = "T"
$ws.Range("B22").Value = "F"
$ws.Range("B23").Value = "S"
$ws.Range("B24").Value = "M"
$ws.Range("B25").Value = "T"
$ws.Range("B26").Value = "W"
$ws.Range("B27").Value = "T"
$ws.Range("B28").Value = "F"
$ws.Range("B29").Value = "M"

# New row for the latest day worked (Tue June 4 2024), added after the insert
# so it lands in the now-shifted columns (A=Date, B=Day, C=Hours, E=Notes).
$ws.Range("A30").Value = 45447
$ws.Range("A30").NumberFormat = "d-mmm"
$ws.Range("B30").NumberFormat = "d-mmm"
$ws.Range("B30").Value = "T"
$ws.Range("C30").Value = 8
$ws.Range("E30").Value = "Sorting code for meeting, meeting, grid search of parameters, added 95% CI error bars to predictions graph"
$ws.Range("E30").WrapText = $true

# Update selection / view to match the author's final state
$ws.Range("E30").Select()

Write-Host "edit complete"
